$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044777178068919
$ws.Range("D2").Value = 1.047077408933838
$ws.Range("E2").Value = 1.052788740798424
$ws.Range("F2").Value = 1.064106469166505
$ws.Range("I2").Value = 1.043271090259124
$ws.Range("J2").Value = 1.049840447941362
$ws.Range("K2").Value = 1.049841055284195
$ws.Range("L2").Value = 1.055536503558201
$ws.Range("M2").Value = 1.06682330188084
$ws.Range("N2").Value = 1.020440124633794
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045775770204787
$ws.Range("D3").Value = 1.047844407548638
$ws.Range("E3").Value = 1.053707367659263
$ws.Range("F3").Value = 1.065177441659183
$ws.Range("I3").Value = 1.04354554601472
$ws.Range("J3").Value = 1.050486190474443
$ws.Range("K3").Value = 1.050419941599203
$ws.Range("L3").Value = 1.056267774767188
$ws.Range("M3").Value = 1.067708773877698
$ws.Range("N3").Value = 1.020659050960056
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04642222776657
$ws.Range("D4").Value = 1.048340914236005
$ws.Range("E4").Value = 1.054302435142328
$ws.Range("F4").Value = 1.065871365662233
$ws.Range("I4").Value = 1.043722126586213
$ws.Range("J4").Value = 1.050903720152407
$ws.Range("K4").Value = 1.050794065477264
$ws.Range("L4").Value = 1.056740979271857
$ws.Range("M4").Value = 1.068282064022743
$ws.Range("N4").Value = 1.020800502503457
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046694070278302
$ws.Range("D5").Value = 1.048549694304685
$ws.Range("E5").Value = 1.054552757332667
$ws.Range("F5").Value = 1.066163313963843
$ws.Range("I5").Value = 1.043796118739348
$ws.Range("J5").Value = 1.051079175060967
$ws.Range("K5").Value = 1.050951237744064
$ws.Range("L5").Value = 1.056939919309736
$ws.Range("M5").Value = 1.068523153634891
$ws.Range("N5").Value = 1.020859918535394
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046739718026292
$ws.Range("D6").Value = 1.048584752222013
$ws.Range("E6").Value = 1.054594796642263
$ws.Range("F6").Value = 1.066212346402792
$ws.Range("I6").Value = 1.043808528127136
$ws.Range("J6").Value = 1.051108630321719
$ws.Range("K6").Value = 1.05097762122426
$ws.Range("L6").Value = 1.056973322483173
$ws.Range("M6").Value = 1.068563638220832
$ws.Range("N6").Value = 1.020869891796355
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046425859858149
$ws.Range("D7").Value = 1.048343703775523
$ws.Range("E7").Value = 1.054305779348493
$ws.Range("F7").Value = 1.065875265815798
$ws.Range("I7").Value = 1.043723116225585
$ws.Range("J7").Value = 1.050906064883243
$ws.Range("K7").Value = 1.050796166051515
$ws.Range("L7").Value = 1.056743637500105
$ws.Range("M7").Value = 1.068285285166611
$ws.Range("N7").Value = 1.020801296620995
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04511459420929
$ws.Range("D8").Value = 1.047336576201693
$ws.Range("E8").Value = 1.053099059013926
$ws.Range("F8").Value = 1.064468215656618
$ws.Range("I8").Value = 1.043364052908581
$ws.Range("J8").Value = 1.050058743249362
$ws.Range("K8").Value = 1.050036786290401
$ws.Range("L8").Value = 1.055783634891148
$ws.Range("M8").Value = 1.067122482422934
$ws.Range("N8").Value = 1.020514154787788
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042806307643039
$ws.Range("D9").Value = 1.045563519076768
$ws.Range("E9").Value = 1.050977714106148
$ws.Range("F9").Value = 1.061995993879366
$ws.Range("I9").Value = 1.042723613220694
$ws.Range("J9").Value = 1.048563315224039
$ws.Range("K9").Value = 1.048695211818713
$ws.Range("L9").Value = 1.054092194483914
$ws.Range("M9").Value = 1.065076037878929
$ws.Range("N9").Value = 1.020006591009487
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041269042832758
$ws.Range("D10").Value = 1.044382633558573
$ws.Range("E10").Value = 1.049566929614703
$ws.Range("F10").Value = 1.060352716803427
$ws.Range("I10").Value = 1.042291481202697
$ws.Range("J10").Value = 1.047564825585463
$ws.Range("K10").Value = 1.04779854691728
$ws.Range("L10").Value = 1.052964745302572
$ws.Range("M10").Value = 1.063713504356682
$ws.Range("N10").Value = 1.019667168012467
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040603772729712
$ws.Range("D11").Value = 1.043871581988783
$ws.Range("E11").Value = 1.048956871598101
$ws.Range("F11").Value = 1.059642323477747
$ws.Range("I11").Value = 1.042103140704174
$ws.Range("J11").Value = 1.047132110808541
$ws.Range("K11").Value = 1.047409747492412
$ws.Range("L11").Value = 1.052476597244612
$ws.Range("M11").Value = 1.063123938341063
$ws.Range("N11").Value = 1.019519949248929
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040356718530626
$ws.Range("D12").Value = 1.043681797578531
$ws.Range("E12").Value = 1.048730392984587
$ws.Range("F12").Value = 1.059378626134811
$ws.Range("I12").Value = 1.042032998989771
$ws.Range("J12").Value = 1.046971327363984
$ws.Range("K12").Value = 1.047265249866517
$ws.Range("L12").Value = 1.052295284629156
$ws.Range("M12").Value = 1.062905010945282
$ws.Range("N12").Value = 1.019465228894726
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040409709922687
$ws.Range("D13").Value = 1.043722505023397
$ws.Range("E13").Value = 1.048778967788221
$ws.Range("F13").Value = 1.059435182209496
$ws.Range("I13").Value = 1.042048052932147
$ws.Range("J13").Value = 1.047005818390066
$ws.Range("K13").Value = 1.047296248718344
$ws.Range("L13").Value = 1.052334176460837
$ws.Range("M13").Value = 1.06295196872426
$ws.Range("N13").Value = 1.019476968255266
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040583349993386
$ws.Range("D14").Value = 1.043855893470028
$ws.Range("E14").Value = 1.048938148269352
$ws.Range("F14").Value = 1.059620522612461
$ws.Range("I14").Value = 1.04209734651321
$ws.Range("J14").Value = 1.047118821495682
$ws.Range("K14").Value = 1.04739780490655
$ws.Range("L14").Value = 1.052461609734893
$ws.Range("M14").Value = 1.063105840431991
$ws.Range("N14").Value = 1.019515426795556
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040690342899037
$ws.Range("D15").Value = 1.043938084205102
$ws.Range("E15").Value = 1.049036241084011
$ws.Range("F15").Value = 1.059734740075833
$ws.Range("I15").Value = 1.042127693586629
$ws.Range("J15").Value = 1.047188439310146
$ws.Range("K15").Value = 1.047460366414333
$ws.Range("L15").Value = 1.052540126566229
$ws.Range("M15").Value = 1.063200654324042
$ws.Range("N15").Value = 1.019539117509599
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041313202559022
$ws.Range("D16").Value = 1.044416556335442
$ws.Range("E16").Value = 1.049607434517459
$ws.Range("F16").Value = 1.060399887687492
$ws.Range("I16").Value = 1.042303954968166
$ws.Range("J16").Value = 1.047593535852357
$ws.Range("K16").Value = 1.047824338975732
$ws.Range("L16").Value = 1.052997143075247
$ws.Range("M16").Value = 1.063752640809079
$ws.Range("N16").Value = 1.019676933268763
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041704006759931
$ws.Range("D17").Value = 1.044716764755645
$ws.Range("E17").Value = 1.049965949398295
$ws.Range("F17").Value = 1.060817427272549
$ws.Range("I17").Value = 1.042414191456232
$ws.Range("J17").Value = 1.047847545600622
$ws.Range("K17").Value = 1.048052505810628
$ws.Range("L17").Value = 1.053283829870738
$ws.Range("M17").Value = 1.064099000195406
$ws.Range("N17").Value = 1.019763315680029
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041931992468514
$ws.Range("D18").Value = 1.044891898171196
$ws.Range("E18").Value = 1.050175144224485
$ws.Range("F18").Value = 1.06106108259132
$ws.Range("I18").Value = 1.042478372369716
$ws.Range("J18").Value = 1.047995670210848
$ws.Range("K18").Value = 1.04818553975629
$ws.Range("L18").Value = 1.053451053720153
$ws.Range("M18").Value = 1.06430106624985
$ws.Range("N18").Value = 1.019813677290644
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042009735839874
$ws.Range("D19").Value = 1.044951618679831
$ws.Range("E19").Value = 1.050246487658441
$ws.Range("F19").Value = 1.061144181659275
$ws.Range("I19").Value = 1.042500236350799
$ws.Range("J19").Value = 1.048046170926843
$ws.Range("K19").Value = 1.048230892072178
$ws.Range("L19").Value = 1.05350807344217
$ws.Range("M19").Value = 1.064369972401999
$ws.Range("N19").Value = 1.019830845255472
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04166207341807
$ws.Range("D20").Value = 1.044684552442096
$ws.Range("E20").Value = 1.049927475950961
$ws.Range("F20").Value = 1.060772617684679
$ws.Range("I20").Value = 1.042402376344516
$ws.Range("J20").Value = 1.047820296372864
$ws.Range("K20").Value = 1.048028031033679
$ws.Range("L20").Value = 1.053253070642689
$ws.Range("M20").Value = 1.06406183491948
$ws.Range("N20").Value = 1.019754050123445
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040532215751792
$ws.Range("D21").Value = 1.043816612705846
$ws.Range("E21").Value = 1.048891270153385
$ws.Range("F21").Value = 1.059565939664861
$ws.Range("I21").Value = 1.042082835849914
$ws.Range("J21").Value = 1.047085546381729
$ws.Range("K21").Value = 1.047367901349968
$ws.Range("L21").Value = 1.052424083605579
$ws.Range("M21").Value = 1.063060527286904
$ws.Range("N21").Value = 1.019504102726374
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039822157812855
$ws.Range("D22").Value = 1.043271153173484
$ws.Range("E22").Value = 1.048240484801778
$ws.Range("F22").Value = 1.058808262452699
$ws.Range("I22").Value = 1.041880865564805
$ws.Range("J22").Value = 1.046623268245719
$ws.Range("K22").Value = 1.04695238783669
$ws.Range("L22").Value = 1.051902909117254
$ws.Range("M22").Value = 1.062431334115035
$ws.Range("N22").Value = 1.019346738380818
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040198541731321
$ws.Range("D23").Value = 1.043560287730028
$ws.Range("E23").Value = 1.048585410085287
$ws.Range("F23").Value = 1.05920982561495
$ws.Range("I23").Value = 1.04198803444645
$ws.Range("J23").Value = 1.046868360004853
$ws.Range("K23").Value = 1.04717270318505
$ws.Range("L23").Value = 1.05217918935025
$ws.Range("M23").Value = 1.062764846181038
$ws.Range("N23").Value = 1.019430180252112
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041681021178245
$ws.Range("D24").Value = 1.044699107716893
$ws.Range("E24").Value = 1.04994486020117
$ws.Range("F24").Value = 1.060792864863971
$ws.Range("I24").Value = 1.042407715450034
$ws.Range("J24").Value = 1.047832609234034
$ws.Range("K24").Value = 1.04803909029228
$ws.Range("L24").Value = 1.053266969400571
$ws.Range("M24").Value = 1.064078628182812
$ws.Range("N24").Value = 1.019758236903149
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04340277618526
$ws.Range("D25").Value = 1.046021697827877
$ws.Range("E25").Value = 1.051525528981329
$ws.Range("F25").Value = 1.062634267129997
$ws.Range("I25").Value = 1.042890095427061
$ws.Range("J25").Value = 1.048950192802465
$ws.Range("K25").Value = 1.049042445960291
$ws.Range("L25").Value = 1.054529443882814
$ws.Range("M25").Value = 1.065604785074316
$ws.Range("N25").Value = 1.020137994245649

Write-Host "Applied 380 kV case updates"